# Add a new weekly batch of "Zapallo" price records for Camote/Paine
# (Mercado Mayorista Lo Valledor de Santiago) by inserting 6 new rows
# right above the current row 1221, pushing the existing data down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows before row 1221 (existing rows 1221-1309 shift to 1227-1315)
$ws.Rows("1221:1226").Insert()

# New data to populate into the freshly inserted rows 1221-1226
$data = @(
    @(44578, "Camote", "1a nueva(o)", 2570, 400, 430, 415, "Región Metropolitana", 415),
    @(44578, "Camote", "1a nueva(o)", 1630, 400, 450, 423, "Región de O'Higgins", 423),
    @(44578, "Camote", "2a nueva(o)", 780,  350, 350, 350, "Región Metropolitana", 350),
    @(44578, "Camote", "2a nueva(o)", 560,  370, 370, 370, "Región de O'Higgins", 370),
    @(44578, "Paine",  "1a nueva(o)", 2300, 150, 180, 166, "Región de O'Higgins", 166),
    @(44578, "Paine",  "2a nueva(o)", 820,  120, 120, 120, "Región de O'Higgins", 120)
)

$r = 1221
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value  = 6
    $ws.Cells.Item($r, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value  = "Metropolitana"
    $ws.Cells.Item($r, 4).Value  = $row[0]
    $ws.Cells.Item($r, 5).Value  = 13
    $ws.Cells.Item($r, 6).Value  = 100112045
    $ws.Cells.Item($r, 7).Value  = "Zapallo"
    $ws.Cells.Item($r, 8).Value  = $row[1]
    $ws.Cells.Item($r, 9).Value  = $row[2]
    $ws.Cells.Item($r, 10).Value = $row[3]
    $ws.Cells.Item($r, 11).Value = $row[4]
    $ws.Cells.Item($r, 12).Value = $row[5]
    $ws.Cells.Item($r, 13).Value = $row[6]
    $ws.Cells.Item($r, 14).Value = "`$/kilo (volumen en unidades)"
    $ws.Cells.Item($r, 15).Value = $row[7]
    $ws.Cells.Item($r, 16).Value = $row[8]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
    $r = $r + 1
}
